# QS_Werte sheet: add "py model und reibung" columns (fvFek, Furnierebene,
# rhok_Fe, G0mean, G4545), rename U column header, and refresh the
# friction-adjusted Iz/A/t values in the lower table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QS_Werte")

# --- Extend the 10 new trailing columns (AB:AK) with the same 15-wide
#     custom width used by every other column on the sheet ----------------
for ($col = 28; $col -le 37; $col++) {
    $ws.Columns.Item($col).ColumnWidth = 14.17
}

# --- Row 1 header row ------------------------------------------------------
# Columns A:H keep their text/position. I1 is a brand-new header inserted
# before the old rhok/E0mean/E90mean block (which shifts right by one, to
# J1:L1), and four more brand-new headers (M1:P1) are inserted before the
# trailing alpha/alphaT pair (which shifts right to Q1:R1).

# New header cell I1 needs the bold/bordered header style - borrow it from
# the existing H1 header cell, then overwrite the text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "fvFek"

# Old J1:L1 (rhok, E0mean, E90mean) already have the header style; just
# rewrite their text in the new shifted-right positions.
$ws.Range("J1").Value = "rhok"
$ws.Range("K1").Value = "E0mean"
$ws.Range("L1").Value = "E90mean"

# Four brand-new headers, again cloning the header style first.
$ws.Range("H1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "Furnierebene"
$ws.Range("H1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "rhok_Fe"
$ws.Range("H1").Copy($ws.Range("O1"))
$ws.Range("O1").Value = "G0mean"
$ws.Range("H1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "G4545"

$ws.Range("H1").Copy($ws.Range("Q1"))
$ws.Range("Q1").Value = "alpha"
$ws.Range("H1").Copy($ws.Range("R1"))
$ws.Range("R1").Value = "alphaT"

# --- Row 2 data row ---------------------------------------------------------
$ws.Range("I2").Value = 31
$ws.Range("J2").Value = 460
$ws.Range("K2").Value = 12000000000
$ws.Range("L2").Value = 370000000
$ws.Range("M2").Value = $true
$ws.Range("N2").Value = 735
$ws.Range("O2").Value = 500
$ws.Range("P2").Value = 4227
$ws.Range("Q2").Value = 0.02
$ws.Range("R2").Value = 0.000004

# --- Table header rename: "U [m]" -> "U_achse [m]" -------------------------
$ws.Range("F6").Value = "U_achse [m]"

# --- Refresh friction-adjusted values (Iz, A, t) for rows 8-17 -------------
$ws.Range("D8").Value = 186.27
$ws.Range("E8").Value = 12.3
$ws.Range("G8").Value = 0.36

$ws.Range("D9").Value = 146.61
$ws.Range("E9").Value = 11.36
$ws.Range("G9").Value = 0.36

$ws.Range("D10").Value = 113.02
$ws.Range("E10").Value = 10.41
$ws.Range("G10").Value = 0.36

$ws.Range("D11").Value = 85
$ws.Range("E11").Value = 9.47
$ws.Range("G11").Value = 0.36

$ws.Range("D12").Value = 62.04
$ws.Range("E12").Value = 8.52
$ws.Range("G12").Value = 0.36

$ws.Range("D13").Value = 43.65
$ws.Range("E13").Value = 7.58
$ws.Range("G13").Value = 0.36

$ws.Range("D14").Value = 29.31
$ws.Range("E14").Value = 6.64
$ws.Range("G14").Value = 0.36

$ws.Range("D15").Value = 18.51
$ws.Range("E15").Value = 5.69
$ws.Range("G15").Value = 0.36

$ws.Range("D16").Value = 10.77
$ws.Range("E16").Value = 4.75
$ws.Range("G16").Value = 0.36

$ws.Range("D17").Value = 5.55
$ws.Range("E17").Value = 3.8
$ws.Range("G17").Value = 0.36
